$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the district names (shared strings content change)
$ws.Range("A1").Value = "Свердловский"
$ws.Range("A2").Value = "Октябрьский"

# Update the numeric values in column B
$ws.Range("B1").Value = 5
$ws.Range("B2").Value = 2

# Move the active selection from A2 to B2
$ws.Range("B2").Select()
